# Insert a new data row for "Zapallo" (Camote, 1a (cosecha)) at row 470,
# pushing the existing rows 470:559 down to 471:560.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(470).Insert()

$ws.Cells.Item(470, 1).Value = 10
$ws.Cells.Item(470, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(470, 3).Value = "La Araucanía"
$ws.Cells.Item(470, 4).Value = 44694
$ws.Cells.Item(470, 5).Value = 9
$ws.Cells.Item(470, 6).Value = 100112045
$ws.Cells.Item(470, 7).Value = "Zapallo"
$ws.Cells.Item(470, 8).Value = "Camote"
$ws.Cells.Item(470, 9).Value = "1a (cosecha)"
$ws.Cells.Item(470, 10).Value = 350
$ws.Cells.Item(470, 11).Value = 500
$ws.Cells.Item(470, 12).Value = 500
$ws.Cells.Item(470, 13).Value = 500
$ws.Cells.Item(470, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(470, 15).Value = "Región del Maule"
$ws.Cells.Item(470, 16).Value = 500
$ws.Cells.Item(470, 17).Value = 1
$ws.Cells.Item(470, 18).Value = "Hortaliza"
